# Refresh the cryptos price/volume table (and two coin-row swaps) per the
# latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel reinterpret strings
# like "1.00" / "0.620" / "39.575.24" as numbers), without leaving any
# residual cell-style change behind once it's done.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "39.575.24"
Set-TextValue "E2" "  +2.01%  "

# Row 3
Set-TextValue "D3" "2.154.06"
Set-TextValue "E3" "  +2.45%  "

# Row 4
Set-TextValue "E4" "  -0.02%  "

# Row 5
Set-TextValue "D5" "226.74"
Set-TextValue "E5" "  -0.30%  "

# Row 6
Set-TextValue "D6" "0.620"
Set-TextValue "E6" "  +0.79%  "

# Row 7
Set-TextValue "D7" "62.70"
Set-TextValue "E7" "  +1.44%  "

# Row 8
Set-TextValue "E8" "  +0.02%  "

# Row 9
Set-TextValue "D9" "0.389"
Set-TextValue "E9" "  +0.19%  "

# Row 10
Set-TextValue "D10" "0.0841"
Set-TextValue "E10" "  -0.08%  "

# Row 11
Set-TextValue "E11" "  -0.12%  "

# Row 12
Set-TextValue "D12" "15.83"
Set-TextValue "E12" "  +0.50%  "

# Row 13
Set-TextValue "D13" "2.477.00"
Set-TextValue "E13" "  +2.67%  "

# Row 14
Set-TextValue "D14" "21.69"
Set-TextValue "E14" "  -1.41%  "

# Row 15
Set-TextValue "D15" "0.803"
Set-TextValue "E15" "  +0.24%  "

# Row 16
Set-TextValue "E16" "  +0.11%  "

# Row 17
Set-TextValue "D17" "2.151.55"
Set-TextValue "E17" "  +3.20%  "

# Row 18
Set-TextValue "D18" "39.536.02"
Set-TextValue "E18" "  +1.80%  "

# Row 19
Set-TextValue "D19" "71.60"
Set-TextValue "E19" "  +0.00%  "

# Row 20
Set-TextValue "D20" "6.04"
Set-TextValue "E20" "  +0.18%  "

# Row 21
Set-TextValue "D21" "0.0₃0848"

# Row 22
Set-TextValue "D22" "227.26"
Set-TextValue "E22" "  +0.17%  "

# Row 23
Set-TextValue "E23" "  +0.01%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D24" "2.32"
Set-TextValue "E24" "  -0.16%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "2.34"
Set-TextValue "E25" "  +1.22%  "

# Row 26
Set-TextValue "D26" "170.61"
Set-TextValue "E26" "  +0.27%  "

# Row 27
Set-TextValue "D27" "9.40"
Set-TextValue "E27" "  -2.51%  "

# Row 28
Set-TextValue "D28" "0.138"
Set-TextValue "E28" "  +1.69%  "

# Row 29
Set-TextValue "D29" "1.43"
Set-TextValue "E29" "  +1.24%  "

# Row 30
Set-TextValue "D30" "19.58"
Set-TextValue "E30" "  +1.17%  "

# Row 31
Set-TextValue "E31" "  +4.39%  "

# Row 32
Set-TextValue "E32" "  +0.81%  "

# Row 33
Set-TextValue "D33" "4.57"
Set-TextValue "E33" "  +0.21%  "

# Row 34
Set-TextValue "D34" "4.69"
Set-TextValue "E34" "  -2.07%  "

# Row 35
Set-TextValue "D35" "6.94"
Set-TextValue "E35" "  -3.04%  "

# Row 36
Set-TextValue "D36" "0.0615"
Set-TextValue "E36" "  +0.01%  "

# Row 37
Set-TextValue "D37" "3.76"
Set-TextValue "E37" "  +7.58%  "

# Row 38
Set-TextValue "D38" "2.38"
Set-TextValue "E38" "  +1.29%  "

# Row 39
Set-TextValue "E39" "  -0.11%  "

# Row 40
Set-TextValue "D40" "4.82"
Set-TextValue "E40" "  +16.16%  "

# Row 41
Set-TextValue "D41" "102.82"
Set-TextValue "E41" "  +1.01%  "

# Row 42
Set-TextValue "E42" "  -1.27%  "

# Row 43
Set-TextValue "D43" "17.61"
Set-TextValue "E43" "  -2.10%  "

# Row 44
Set-TextValue "D44" "1.513.15"
Set-TextValue "E44" "  -0.75%  "

# Row 45
Set-TextValue "E45" "  -0.38%  "

# Row 46
Set-TextValue "D46" "7.86"
Set-TextValue "E46" "  +1.09%  "

# Row 47
Set-TextValue "D47" "2.80"
Set-TextValue "E47" "  -0.05%  "

# Row 48
Set-TextValue "D48" "0.0919"
Set-TextValue "E48" "  +0.92%  "

# Row 49
Set-TextValue "E49" "  +0.40%  "

# Row 50
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue "D50" "0.000191"
Set-TextValue "E50" "  +28.72%  "

# Row 51
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D51" "2.98"
Set-TextValue "E51" "  +1.07%  "
